$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.326.74"
$ws.Range("E2").Value = "  -3.28%  "
$ws.Range("D3").Value = "3.492.66"
$ws.Range("E3").Value = "  -4.83%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'606.21"
$ws.Range("E5").Value = "  -2.26%  "
$ws.Range("D6").Value = "'148.54"
$ws.Range("E6").Value = "  -6.97%  "
$ws.Range("D7").Value = "3.492.21"
$ws.Range("E7").Value = "  -4.77%  "
$ws.Range("D9").Value = "'0.480"
$ws.Range("E9").Value = "  -3.14%  "
$ws.Range("E10").Value = "  -3.85%  "
$ws.Range("D11").Value = "'6.97"
$ws.Range("E11").Value = "  -3.12%  "
$ws.Range("E12").Value = "  -3.93%  "
$ws.Range("D13").Value = "'0.0000218"
$ws.Range("E13").Value = "  -4.70%  "
$ws.Range("D14").Value = "4.082.59"
$ws.Range("E14").Value = "  -4.80%  "
$ws.Range("D15").Value = "'31.43"
$ws.Range("E15").Value = "  -2.90%  "
$ws.Range("D16").Value = "3.494.70"
$ws.Range("E16").Value = "  -5.89%  "
$ws.Range("D17").Value = "67.220.24"
$ws.Range("E17").Value = "  -3.50%  "
$ws.Range("D18").Value = "'0.117"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("D19").Value = "'6.38"
$ws.Range("E19").Value = "  -2.00%  "
$ws.Range("D20").Value = "'15.06"
$ws.Range("E20").Value = "  -4.83%  "
$ws.Range("D21").Value = "'446.89"
$ws.Range("E21").Value = "  -5.09%  "
$ws.Range("E22").Value = "  -12.75%  "
$ws.Range("D23").Value = "'0.621"
$ws.Range("E23").Value = "  -4.19%  "
$ws.Range("D24").Value = "'77.08"
$ws.Range("E24").Value = "  -3.47%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "'0.0000129"
$ws.Range("E25").Value = "  +4.42%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").Value = "3.632.03"
$ws.Range("E27").Value = "  -4.85%  "
$ws.Range("E28").Value = "  -8.42%  "
$ws.Range("D29").Value = "'8.29"
$ws.Range("E29").Value = "  -4.64%  "
$ws.Range("D30").Value = "'2.47"
$ws.Range("E30").Value = "  -4.49%  "
$ws.Range("D31").Value = "'1.56"
$ws.Range("E31").Value = "  -6.43%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "'0.162"
$ws.Range("E33").Value = "  -0.76%  "
$ws.Range("D34").Value = "'25.66"
$ws.Range("D35").Value = "'6.14"
$ws.Range("E35").Value = "  -3.64%  "
$ws.Range("E36").Value = "  -6.30%  "
$ws.Range("D37").Value = "3.481.54"
$ws.Range("E37").Value = "  -5.19%  "
$ws.Range("D38").Value = "'8.00"
$ws.Range("E38").Value = "  -3.10%  "
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").Value = "'2.19"
$ws.Range("E41").Value = "  -0.90%  "
$ws.Range("D42").Value = "'173.51"
$ws.Range("E42").Value = "  -3.33%  "
$ws.Range("D43").Value = "'0.0875"
$ws.Range("E43").Value = "  -1.51%  "
$ws.Range("E44").Value = "  -6.62%  "
$ws.Range("D45").Value = "'0.882"
$ws.Range("E45").Value = "  -4.65%  "
$ws.Range("E46").Value = "  -2.84%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'27.30"
$ws.Range("E47").Value = "  -5.13%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "'1.26"
$ws.Range("E48").Value = "  +4.19%  "
$ws.Range("D49").Value = "'2.56"
$ws.Range("E49").Value = "  -6.05%  "
$ws.Range("D50").Value = "'7.54"
$ws.Range("E50").Value = "  -3.77%  "
$ws.Range("E51").Value = "  -2.99%  "
